# Adiciona a coluna "País" e mais linhas de dados (Nome, Idade, Pais)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dados completos da tabela (incluindo cabecalho) - cada linha = (Nome, Idade, Pais)
$dados = @(
    @("Nome", "Idade", "País"),
    @("Will", 18, "Brasil"),
    @("ChingChongBingBong", 29, "China"),
    @("Kyle", 8, "US")
)

$linhaInicial = 1
$numLinhas = $dados.Count
$numColunas = 3

for ($i = 0; $i -lt $numLinhas; $i++) {
    $linha = $linhaInicial + $i
    for ($j = 0; $j -lt $numColunas; $j++) {
        $coluna = $j + 1
        $valor = $dados[$i][$j]
        if ($valor -eq $null) {
            # pula celulas vazias
            continue
        }
        $ws.Cells.Item($linha, $coluna).Value = $valor
    }
}

$ultimaLinha = $linhaInicial + $numLinhas - 1
$ws.Cells.Item($ultimaLinha, $numColunas).Select()
